$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 26 (current "877 / 청계현대" row),
# shifting all following rows down by 2.
$ws.Rows.Item(26).Resize(2).Insert()

# New row 26: 950 / 행당대림
$ws.Cells.Item(26, 1).Value = 950
$ws.Cells.Item(26, 2).Value = "행당대림"

# New row 27: 878 / 행당한진
$ws.Cells.Item(27, 1).Value = 878
$ws.Cells.Item(27, 2).Value = "행당한진"

# Match the style/font used for similar manually-added rows (e.g. row 21)
# and the row height seen in the diff (13.2).
$ws.Range("B26:B27").Font.Name = "맑은 고딕"
$ws.Rows.Item(26).RowHeight = 13.2
$ws.Rows.Item(27).RowHeight = 13.2

# Update the sheet view: drop the frozen/scrolled topLeftCell and move
# the active selection.
$ws.Range("G23").Select()
